# Update the "dSF" (column F) values to reflect repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -3
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = -8
$ws.Range("F10").Value = -8
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = -3
$ws.Range("F15").Value = 0
$ws.Range("F18").Value = -1
